$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($range, $value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-CellText "D2" "36.672.82"
Set-CellText "E2" "  -0.32%  "
Set-CellText "D3" "2.114.17"
Set-CellText "E3" "  +9.65%  "
Set-CellText "E4" "  +0.04%  "
Set-CellText "D5" "254.04"
Set-CellText "E5" "  +1.34%  "
Set-CellText "D6" "0.665"
Set-CellText "E6" "  -5.50%  "
Set-CellText "E7" "  +0.03%  "
Set-CellText "D8" "47.72"
Set-CellText "E8" "  +7.62%  "
Set-CellText "D9" "60.43"
Set-CellText "E9" "  +2.81%  "
Set-CellText "D10" "0.374"
Set-CellText "E10" "  +1.00%  "
Set-CellText "D11" "0.0747"
Set-CellText "E11" "  -2.80%  "
Set-CellText "E12" "  +0.06%  "
Set-CellText "D13" "2.419.54"
Set-CellText "E13" "  +9.64%  "
Set-CellText "D14" "14.32"
Set-CellText "E14" "  -2.99%  "
Set-CellText "D15" "0.833"
Set-CellText "E15" "  +1.28%  "
Set-CellText "D16" "2.109.79"
Set-CellText "E16" "  +9.29%  "
Set-CellText "E17" "  -0.52%  "
Set-CellText "D18" "36.695.92"
Set-CellText "E18" "  +0.00%  "
Set-CellText "D19" "73.59"
Set-CellText "E19" "  -1.39%  "
Set-CellText "E20" "  -3.69%  "
Set-CellText "D21" "13.26"
Set-CellText "E21" "  -1.47%  "
Set-CellText "D22" "240.71"
Set-CellText "E22" "  -4.62%  "
Set-CellText "D23" "5.20"
Set-CellText "E23" "  -0.92%  "
Set-CellText "E24" "  +0.02%  "
Set-CellText "E25" "  -7.56%  "
Set-CellText "D26" "172.84"
Set-CellText "E26" "  +2.67%  "
Set-CellText "D27" "21.60"
Set-CellText "E27" "  +14.68%  "
Set-CellText "E28" "  +3.23%  "
Set-CellText "E29" "  -9.55%  "
Set-CellText "D30" "29.78"
Set-CellText "E30" "  +66.57%  "
Set-CellText "D31" "0.124"
Set-CellText "E31" "  -4.81%  "
Set-CellText "D32" "4.50"
Set-CellText "E32" "  -1.88%  "
Set-CellText "D33" "0.0603"
Set-CellText "E33" "  -2.64%  "
Set-CellText "D34" "0.0918"
Set-CellText "E34" "  +5.55%  "
Set-CellText "D35" "0.965"
Set-CellText "E35" "  +7.26%  "
Set-CellText "D36" "2.36"
Set-CellText "E36" "  +14.46%  "
Set-CellText "E37" "  -4.22%  "
Set-CellText "E39" "  -6.53%  "
Set-CellText "E40" "  -12.32%  "
Set-CellText "E41" "  +6.44%  "
Set-CellText "E42" "  -1.71%  "
Set-CellText "D43" "98.77"
Set-CellText "E43" "  -8.21%  "
Set-CellText "E44" "  +7.27%  "
Set-CellText "D45" "16.02"
Set-CellText "E45" "  -7.91%  "
Set-CellText "D46" "1.347.93"
Set-CellText "E46" "  +0.18%  "
Set-CellText "B47" "FraxShare"
Set-CellText "C47" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-CellText "D47" "7.20"
Set-CellText "E47" "  +11.52%  "
Set-CellText "B48" "Cronos"
Set-CellText "C48" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-CellText "D48" "0.0844"
Set-CellText "E48" "  +3.30%  "
Set-CellText "D49" "2.302.48"
Set-CellText "E50" "  +0.99%  "
Set-CellText "E51" "  -4.96%  "
